$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)
# D-column values are forced to remain text (avoiding Excel's auto numeric coercion)
# by temporarily applying a text NumberFormat, then clearing formatting back to default.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.910.09'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.635.92'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.59'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5068'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06372'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07751'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.274'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.620.50'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5539'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.64%  '
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.18'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.915.23'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.445'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '195.36'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.913'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.065'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.902'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.77'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1245'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +6.08%  '
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('E28').Value = '  -0.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.245'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04875'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.255'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('E32').Value = '  +0.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.546'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.371'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9066'
$ws.Range('D35').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.568'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5507'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.121.57'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01560'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.001'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.583'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '97.81'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.84%  '
$ws.Range('E44').Value = '  -3.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.773.04'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4455'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.93'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9972'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05145'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.538'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.005'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.22%  '
